$d = $word.ActiveDocument

# Locate the two empty paragraphs right after the "Enemies:" heading by
# scanning the paragraph collection (robust to any off-by-one surprises).
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Enemies:") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'Enemies:' paragraph"
}

$firstEmptyIndex = $targetIndex + 1
$secondEmptyIndex = $targetIndex + 2

# --- Fill in the two new enemy-description paragraphs -----------------
$pJump = $d.Paragraphs.Item($firstEmptyIndex)
$pJump.Range.Text = "Jumping enemy: Walks  across the platforms, occasionally jumps which then causes platforms under the Players (random pick of the 2 in multiplayer) to begin to fall / be destroyed. - Tom"

$pSnow = $d.Paragraphs.Item($secondEmptyIndex)

# Re-anchor the _GoBack bookmark (it previously sat in its own empty
# paragraph right before "AI for Player 2's character") onto this
# paragraph first, while it is still empty, then insert the paragraph's
# text with InsertBefore so the bookmark ends up anchored at the end of
# the new text (matching Word's normal gravity behaviour for bookmarks
# when text is typed in front of them). Re-using the existing bookmark
# name automatically relocates it and removes it from its old spot.
$d.Bookmarks.Add("_GoBack", $pSnow.Range)
$pSnow = $d.Paragraphs.Item($secondEmptyIndex)
$pSnow.Range.InsertBefore("Snowball Enemy: Shoots snowball across the line that will slow  / freeze the player. - Jasper")
